$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the CDS read table function (query text) in B3:
#  - use samp.sample_tumor_status instead of the undefined `tumor` alias
#  - re-indent the trailing ORDER BY clause
$newQuery = @'
MATCH (s:study)<--(p:participant)<--(samp:sample)
WHERE s.study_name in ["GECCO OICR: Molecular Pathological Epidemiology of Colorectal Cancer"]
WITH p,s,samp,COLLECT(DISTINCT samp.sample_tumor_status) as tumor
RETURN  
 coalesce(samp.sample_id, '') as `Sample ID`,
 coalesce(p.participant_id,'') as `Participant ID`,
 coalesce(s.study_name, '') as `Study Name`,
 coalesce(s.phs_accession,'') as `Accession`,
 coalesce(samp.sample_tumor_status,'') as `Tumor`,
coalesce(samp.sample_type,'') as `Analyte Type`
  ORDER By samp.sample_id LIMIT 100
'@

$ws.Range("B3").Value = $newQuery

# Move the active selection from B4 to B3
$ws.Range("B3").Select()
